$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 30.92471370227826
$ws.Range("B3").Value = 30.35570958106607
$ws.Range("B4").Value = 39.9365689147778
$ws.Range("I5").Value = 96.20297504441875
$ws.Range("I6").Value = 96.3169143589794
$ws.Range("I7").Value = 96.33566636993102
$ws.Range("C8").Value = 50.5068007364945
$ws.Range("C9").Value = 50.43083459661063
$ws.Range("C10").Value = 50.64491845601064
$ws.Range("D11").Value = 98.69383719463904
$ws.Range("D12").Value = 98.78253175939945
$ws.Range("D13").Value = 98.97945518055859
$ws.Range("E14").Value = 98.92964338564892
$ws.Range("E15").Value = 98.99039607038232
$ws.Range("E16").Value = 98.93362909456415
$ws.Range("F17").Value = 98.46251930039219
$ws.Range("F18").Value = 98.57204363408384
$ws.Range("F19").Value = 98.44133171789868
$ws.Range("G20").Value = 98.00793807320684
$ws.Range("G21").Value = 98.10782911548975
$ws.Range("G22").Value = 98.04275608301113
$ws.Range("H23").Value = 97.5851356131321
$ws.Range("H24").Value = 97.56652208834824
$ws.Range("H25").Value = 97.5766935032821
$ws.Range("B26").Value = 41.06677166725142
$ws.Range("B27").Value = 47.36611267523067
$ws.Range("I28").Value = 96.16475911118812
$ws.Range("I29").Value = 96.14163516294354
$ws.Range("C30").Value = 50.21696702457879
$ws.Range("C31").Value = 50.37192329591493
$ws.Range("D32").Value = 98.84921127727753
$ws.Range("D33").Value = 98.79820028349364
$ws.Range("E34").Value = 98.95035318381366
$ws.Range("E35").Value = 98.9480933475057
$ws.Range("F36").Value = 98.46046561045164
$ws.Range("F37").Value = 98.56934767549411
$ws.Range("G38").Value = 97.98752707599905
$ws.Range("G39").Value = 98.02470566324116
$ws.Range("H40").Value = 97.62682939546238
$ws.Range("H41").Value = 97.70536607346132
$ws.Range("B42").Value = 36.79501782827266
$ws.Range("B43").Value = 49.84724520074579
$ws.Range("I44").Value = 96.23373656532033
$ws.Range("I45").Value = 96.38856597596097
$ws.Range("C46").Value = 50.76713118922515
$ws.Range("C47").Value = 50.09554112996507
$ws.Range("D48").Value = 98.95920040536885
$ws.Range("D49").Value = 98.79746150758228
$ws.Range("E50").Value = 98.98328389254418
$ws.Range("E51").Value = 98.85107360260405
$ws.Range("F52").Value = 98.55347271192133
$ws.Range("F53").Value = 98.41268829169429
$ws.Range("G54").Value = 97.99794894977717
$ws.Range("G55").Value = 98.04066174835239
$ws.Range("H56").Value = 97.49132686849759
$ws.Range("H57").Value = 97.48134028923215
